$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5280971303759259
$ws.Range("C2").Value = 0.2284544404074147
$ws.Range("D2").Value = 0.06097336337090553
$ws.Range("E2").Value = 0.1388871599331836
$ws.Range("F2").Value = 1.181426164879099
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.26198650538376
$ws.Range("L2").Value = 0.1918557080100527
$ws.Range("M2").Value = 0.154433964690341
$ws.Range("O2").Value = 4.300780462613432
$ws.Range("B3").Value = 0.4935342816661148
$ws.Range("C3").Value = 0.2284731720940627
$ws.Range("D3").Value = 0.05921214971375832
$ws.Range("E3").Value = 0.1394203463867356
$ws.Range("F3").Value = 1.182560946228818
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.2308765021766277
$ws.Range("L3").Value = 0.1893104981600615
$ws.Range("M3").Value = 0.14793915479855
$ws.Range("O3").Value = 4.319256954297458
$ws.Range("B4").Value = 0.4724925182315758
$ws.Range("C4").Value = 0.228500914113404
$ws.Range("D4").Value = 0.05811692019206305
$ws.Range("E4").Value = 0.1397997322016362
$ws.Range("F4").Value = 1.183865847586247
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.2117681803739515
$ws.Range("L4").Value = 0.1878389323733245
$ws.Range("M4").Value = 0.1440178982112492
$ws.Range("O4").Value = 4.3326428784614
$ws.Range("B5").Value = 0.4639636273502106
$ws.Range("C5").Value = 0.2285163272787401
$ws.Range("D5").Value = 0.05766714504269999
$ws.Range("E5").Value = 0.1399674318446973
$ws.Range("F5").Value = 1.184550604753589
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.2039801166259423
$ws.Range("L5").Value = 0.1872622499613357
$ws.Range("M5").Value = 0.1424368054558123
$ws.Range("O5").Value = 4.338611232546157
$ws.Range("B6").Value = 0.4625501937114507
$ws.Range("C6").Value = 0.2285191354574181
$ws.Range("D6").Value = 0.0575922518011609
$ws.Range("E6").Value = 0.1399960698248606
$ws.Range("F6").Value = 1.184673551213663
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.2026868500055201
$ws.Range("L6").Value = 0.1871678827139647
$ws.Range("M6").Value = 0.1421752871704491
$ws.Range("O6").Value = 4.339633291634669
$ws.Range("B7").Value = 0.472377308465866
$ws.Range("C7").Value = 0.2285011053154982
$ws.Range("D7").Value = 0.05811086835243628
$ws.Range("E7").Value = 0.1398019408034727
$ws.Range("F7").Value = 1.183874462871678
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.2116631524149994
$ws.Range("L7").Value = 0.1878310618583683
$ws.Range("M7").Value = 0.1439965066677438
$ws.Range("O7").Value = 4.332721290467617
$ws.Range("B8").Value = 0.5161428231540697
$ws.Range("C8").Value = 0.2284575427020457
$ws.Range("D8").Value = 0.06036897987345213
$ws.Range("E8").Value = 0.1390602199416655
$ws.Range("F8").Value = 1.181691261537352
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.2512613907876471
$ws.Range("L8").Value = 0.1909592244000109
$ws.Range("M8").Value = 0.1521807975338625
$ws.Range("O8").Value = 4.306727593939939
$ws.Range("B9").Value = 0.603375553896683
$ws.Range("C9").Value = 0.2284999591784
$ws.Range("D9").Value = 0.06468667661751226
$ws.Range("E9").Value = 0.1380175735433795
$ws.Range("F9").Value = 1.182233243122447
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.3288466244276265
$ws.Range("L9").Value = 0.1978151420177596
$ws.Range("M9").Value = 0.1687546674975096
$ws.Range("O9").Value = 4.27194527623476
$ws.Range("B10").Value = 0.6683052358698376
$ws.Range("C10").Value = 0.2286077519738399
$ws.Range("D10").Value = 0.06779093721596752
$ws.Range("E10").Value = 0.137501650542287
$ws.Range("F10").Value = 1.185570924565738
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 0.3857946177261908
$ws.Range("L10").Value = 0.2032901012923389
$ws.Range("M10").Value = 0.1812475312340851
$ws.Range("O10").Value = 4.256257972841098
$ws.Range("B11").Value = 0.6980221968194371
$ws.Range("C11").Value = 0.2286731601045133
$ws.Range("D11").Value = 0.06918827726064336
$ws.Range("E11").Value = 0.1373210524912452
$ws.Range("F11").Value = 1.187727542644694
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 0.4116876050145777
$ws.Range("L11").Value = 0.2058755181219141
$ws.Range("M11").Value = 0.1869987861070896
$ws.Range("O11").Value = 4.251263587692108
$ws.Range("B12").Value = 0.7093006842721081
$ws.Range("C12").Value = 0.2287002600376695
$ws.Range("D12").Value = 0.0697152689615379
$ws.Range("E12").Value = 0.1372604262445059
$ws.Range("F12").Value = 1.188635940907702
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 0.4214904134900905
$ws.Range("L12").Value = 0.206868136836647
$ws.Range("M12").Value = 0.1891863525037607
$ws.Range("O12").Value = 4.249680246584859
$ws.Range("B13").Value = 0.7068705433752882
$ws.Range("C13").Value = 0.2286943202625267
$ws.Range("D13").Value = 0.06960186786816536
$ws.Range("E13").Value = 0.1372731382293608
$ws.Range("F13").Value = 1.188436222033957
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 0.4193793127943195
$ws.Range("L13").Value = 0.2066537557146972
$ws.Range("M13").Value = 0.1887147917506127
$ws.Range("O13").Value = 4.250007553553615
$ws.Range("B14").Value = 0.6989495806644754
$ws.Range("C14").Value = 0.2286753430418926
$ws.Range("D14").Value = 0.06923167641241434
$ws.Range("E14").Value = 0.1373159092617193
$ws.Range("F14").Value = 1.187800438933309
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 0.4124941366732457
$ws.Range("L14").Value = 0.2059569096864209
$ws.Range("M14").Value = 0.187178564944027
$ws.Range("O14").Value = 4.251127154418469
$ws.Range("B15").Value = 0.6941010423522869
$ws.Range("C15").Value = 0.2286640218686387
$ws.Range("D15").Value = 0.06900464272985829
$ws.Range("E15").Value = 0.1373431181355009
$ws.Range("F15").Value = 1.187422948328503
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 0.4082764543383348
$ws.Range("L15").Value = 0.205531837807456
$ws.Range("M15").Value = 0.1862388403273059
$ws.Range("O15").Value = 4.251853040963965
$ws.Range("B16").Value = 0.6663667373215105
$ws.Range("C16").Value = 0.2286038048814518
$ws.Range("D16").Value = 0.06769931799706086
$ws.Range("E16").Value = 0.1375145401930524
$ws.Range("F16").Value = 1.1854428241476
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 0.3841021518123
$ws.Range("L16").Value = 0.203123041408162
$ws.Range("M16").Value = 0.1808730355445078
$ws.Range("O16").Value = 4.256627461071218
$ws.Range("B17").Value = 0.6493983586543948
$ws.Range("C17").Value = 0.2285710406288359
$ws.Range("D17").Value = 0.06689473637107568
$ws.Range("E17").Value = 0.1376335458175735
$ws.Range("F17").Value = 1.184391527144953
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 0.3692683701772523
$ws.Range("L17").Value = 0.2016695691995096
$ws.Range("M17").Value = 0.1775986713559163
$ws.Range("O17").Value = 4.260104941671926
$ws.Range("B18").Value = 0.6396555847836964
$ws.Range("C18").Value = 0.2285537379304472
$ws.Range("D18").Value = 0.06643057021173604
$ws.Range("E18").Value = 0.1377070876661985
$ws.Range("F18").Value = 1.183846930933967
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.3607351796930232
$ws.Range("L18").Value = 0.2008424998464449
$ws.Range("M18").Value = 0.1757217687344621
$ws.Range("O18").Value = 4.262306701687379
$ws.Range("B19").Value = 0.6363597845474374
$ws.Range("C19").Value = 0.2285481451565161
$ws.Range("D19").Value = 0.06627317310227454
$ws.Range("E19").Value = 0.1377328629085675
$ws.Range("F19").Value = 1.183672860715404
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.3578457968743862
$ws.Range("L19").Value = 0.2005640038084664
$ws.Range("M19").Value = 0.1750873887549318
$ws.Range("O19").Value = 4.263086807829268
$ws.Range("B20").Value = 0.6512029174816689
$ws.Range("C20").Value = 0.2285743689709037
$ws.Range("D20").Value = 0.06698052969372981
$ws.Range("E20").Value = 0.1376203504858644
$ws.Range("F20").Value = 1.18449722161094
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 0.3708475780983633
$ws.Range("L20").Value = 0.2018233700639627
$ws.Range("M20").Value = 0.1779465688739776
$ws.Range("O20").Value = 4.259713893052719
$ws.Range("B21").Value = 0.7012754768942671
$ws.Range("C21").Value = 0.2286808540230254
$ws.Range("D21").Value = 0.06934046915841208
$ws.Range("E21").Value = 0.1373031358403409
$ws.Range("F21").Value = 1.187984694628099
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 0.4145165458491817
$ws.Range("L21").Value = 0.2061612222386486
$ws.Range("M21").Value = 0.1876295298962916
$ws.Range("O21").Value = 4.250789944102905
$ws.Range("B22").Value = 0.7341480999820362
$ws.Range("C22").Value = 0.2287640262116355
$ws.Range("D22").Value = 0.07087028030135656
$ws.Range("E22").Value = 0.1371410546165386
$ws.Range("F22").Value = 1.190798630028709
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 0.4430431074108867
$ws.Range("L22").Value = 0.2090753718650404
$ws.Range("M22").Value = 0.19401433820601
$ws.Range("O22").Value = 4.246752411474233
$ws.Range("B23").Value = 0.7165900754359029
$ws.Range("C23").Value = 0.2287184008072245
$ws.Range("D23").Value = 0.070054946452629
$ws.Range("E23").Value = 0.137223426817382
$ws.Range("F23").Value = 1.189247874675303
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 0.4278193381064455
$ws.Range("L23").Value = 0.207512816085611
$ws.Range("M23").Value = 0.1906015199495528
$ws.Range("O23").Value = 4.248743119533856
$ws.Range("B24").Value = 0.6503870372293932
$ws.Range("C24").Value = 0.2285728594479224
$ws.Range("D24").Value = 0.06694174751909543
$ws.Range("E24").Value = 0.137626300133693
$ws.Range("F24").Value = 1.184449250838497
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 0.3701336339010481
$ws.Range("L24").Value = 0.201753810053134
$ws.Range("M24").Value = 0.1777892670387047
$ws.Range("O24").Value = 4.259890055314429
$ws.Range("B25").Value = 0.5796279470783077
$ws.Range("C25").Value = 0.2284749341384327
$ws.Range("D25").Value = 0.06353051403530685
$ws.Range("E25").Value = 0.1382556448363506
$ws.Range("F25").Value = 1.181570350968265
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.3078662799450171
$ws.Range("L25").Value = 0.1958833880416435
$ws.Range("M25").Value = 0.1642152292170067
$ws.Range("O25").Value = 4.279621775039516
